# Error Calculations and Plots
# Update the CDF missing-data worksheet: re-impute several cells in columns
# E (D-header) and F, and remove two trailing rows (the data shifted up by
# one row starting at row 26, with row 34/35 no longer present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ---
$ws.Range("F5").ClearContents()

# --- Row 11 ---
$ws.Range("F11").Value = 17.65

# --- Row 19 ---
$ws.Range("E19").Value = -6.5
$ws.Range("F19").ClearContents()

# --- Row 21 ---
$ws.Range("E21").ClearContents()

# --- Row 23 ---
$ws.Range("E23").Value = -7
$ws.Range("F23").Value = 16.48

# --- Row 25 ---
$ws.Range("F25").Value = 16.6

# --- Row 26 (was "RM 232") ---
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

# --- Row 27 (was "SC 5") ---
$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()
$ws.Range("F27").ClearContents()

# --- Row 28 (was "SC 92") ---
$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

# --- Row 29 (was "SC 101") ---
$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = -19.5
$ws.Range("C29").Value = 11.2
$ws.Range("E29").Value = -6.8
$ws.Range("F29").ClearContents()

# --- Row 30 (was "SC 105") ---
$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

# --- Row 31 (was "SC 119") ---
$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

# --- Row 32 (was "SC 120") ---
$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

# --- Row 33 (was "SC 132") ---
$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53

# --- Remove the old trailing rows (SC 193 / SC 232) now that their data has
#     moved up into rows 32/33 above ---
$ws.Range("A35:F35").EntireRow.Delete()
$ws.Range("A34:F34").EntireRow.Delete()
